# Mise à jour résultats avec marge 20%
$wb = $excel.ActiveWorkbook

# --- Sheet "Statistiques" ---
# These cells are stored as text in the workbook (even though some look
# numeric), so a leading apostrophe is used to force text entry and keep
# Excel from auto-converting them to numeric cells.
$ws1 = $wb.Worksheets.Item("Statistiques")

$ws1.Range("B4").Value = "'7001"
$ws1.Range("B5").Value = "'333.02s"
$ws1.Range("B12").Value = "'1154 (99.4%)"
$ws1.Range("B13").Value = "'7 (0.6%)"
$ws1.Range("B15").Value = "'4117"
$ws1.Range("B16").Value = "'527"
$ws1.Range("B17").Value = "'3.55"
$ws1.Range("B18").Value = "'0.45"

# --- Sheet "Par Présentation" ---
$ws2 = $wb.Worksheets.Item("Par Présentation")

$ws2.Range("D2").Value = 4

$ws2.Range("C3").Value = 217
$ws2.Range("D3").Value = 5

$ws2.Range("D5").Value = 5

$ws2.Range("C6").Value = 211
$ws2.Range("D6").Value = 4

$ws2.Range("D7").Value = 4

$ws2.Range("D8").Value = 5

$ws2.Range("D10").Value = 4

$ws2.Range("C13").Value = 209

$ws2.Range("C16").Value = 225

$ws2.Range("D18").Value = 4

$ws2.Range("D19").Value = 2

$ws2.Range("D20").Value = 4

$ws2.Range("D22").Value = 2

$ws2.Range("C23").Value = 52

$ws2.Range("D24").Value = 3

$ws2.Range("D25").Value = 2

$ws2.Range("D26").Value = 3

$ws2.Range("C27").Value = 55

$ws2.Range("C30").Value = 41

$ws2.Range("C32").Value = 43
$ws2.Range("D32").Value = 1
